# Enhance customer template handler to include primary and secondary
# address fields: append 20 new header columns (I1:AB1) after the
# existing H1 header, re-using the same header style (bold, bordered,
# centered) that the other header cells already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell whose formatting (font/border/alignment) the new header
# cells should inherit.
$styleSource = $ws.Range("H1")

$newHeaders = @(
    "customer_primary_address.address_title",
    "customer_primary_address.address_line1",
    "customer_primary_address.address_line2",
    "customer_primary_address.city",
    "customer_primary_address.state",
    "customer_primary_address.zipcode",
    "customer_primary_address.country",
    "customer_primary_address.phone",
    "customer_primary_address.email",
    "customer_primary_address.address_type",
    "customer_secondary_address.address_title",
    "customer_secondary_address.address_line1",
    "customer_secondary_address.address_line2",
    "customer_secondary_address.city",
    "customer_secondary_address.state",
    "customer_secondary_address.zipcode",
    "customer_secondary_address.country",
    "customer_secondary_address.phone",
    "customer_secondary_address.email",
    "customer_secondary_address.address_type"
)

# H1 is column 8, so new headers start at column 9 (I).
$col = 9
foreach ($header in $newHeaders) {
    $cell = $ws.Cells.Item(1, $col)
    $styleSource.Copy($cell)
    $cell.Value2 = $header
    $col = $col + 1
}
